# "Landscaping Data" update: append 7 new observation rows (191-197), all
# dated 2025-06-06 (Excel serial 45814), to the bottom of Sheet1's table,
# and move the view/selection down to the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastOldRow = 190
$firstNewRow = 191
$lastNewRow = 197

# Give the new rows the same cell formatting as the previous last row
# (row 190) before populating them - this is what makes column A (the
# Date column) pick up the existing date number format/style instead of
# defaulting to "General".
$ws.Range("A$lastOldRow`:T$lastOldRow").Copy()
$ws.Range("A$firstNewRow`:T$lastNewRow").PasteSpecial(-4122)

# New row data (columns A-T; F is a formula and is handled separately).
$newRows = @(
    @{ Row=191; A=45814; B="Flowering";     C="Large";  D=65; E=77; G=5.57; H=0.1;  I="Yes"; J=2; K="Dark";    L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=192; A=45814; B="Nonflowering";  C="Medium"; D=65; E=77; G=5.57; H=0.25; I="Yes"; J=3; K="Dark";    L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=193; A=45814; B="Nonflowering";  C="Small";  D=65; E=77; G=5.57; H=0.35; I="Yes"; J=3; K="Dark";    L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=194; A=45814; B="Nonflowering";  C="Medium"; D=65; E=77; G=5.57; H=0.4;  I="Yes"; J=3; K="Dark";    L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=195; A=45814; B="Nonflowering";  C="Medium"; D=65; E=77; G=5.57; H=0.5;  I="Yes"; J=3; K="Neutral"; L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=196; A=45814; B="Nonflowering";  C="Large";  D=65; E=77; G=5.57; H=0.1;  I="Yes"; J=4; K="Neutral"; L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
    @{ Row=197; A=45814; B="Tree";          C="Medium"; D=65; E=77; G=5.57; H=1.25; I="Yes"; J=1; K="Dark";    L=7; M=0.82; N=69; O=29.87; P=9; Q=0.83; R=5.7; S=68; T=47 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Formula = "=ABS(D$row-E$row)"
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

$wb.Application.Calculate()

# Scroll the view down and select N191:N197 (N191 active), matching where
# the author was working after adding the new rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = $lastOldRow
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N$firstNewRow`:N$lastNewRow").Select()
